$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking columns (Price / Volume%) must be forced to Text so they
# round-trip as literal strings (matching the source data export), not as
# auto-coerced Excel numbers/percentages.
function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextCell "D2" "306.26"
Set-TextCell "E2" "-4.79%"
Set-TextCell "D3" "39.84"
Set-TextCell "E3" "-7.12%"
Set-TextCell "D4" "5.082"
Set-TextCell "E4" "-1.79%"
Set-TextCell "D5" "0.07688"
Set-TextCell "E5" "-5.73%"
Set-TextCell "D6" "4.272"
Set-TextCell "E6" "-1.21%"
Set-TextCell "D7" "1.615"
Set-TextCell "E7" "-10.93%"
Set-TextCell "D8" "0.8786"
Set-TextCell "E8" "-6.95%"
Set-TextCell "D9" "0.09648"
Set-TextCell "E9" "-13.68%"
Set-TextCell "D10" "0.1723"
Set-TextCell "E10" "-6.87%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell "D11" "0.08907"
Set-TextCell "E11" "-4.86%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell "D12" "0.04409"
Set-TextCell "E12" "-4.62%"
Set-TextCell "D13" "0.1056"
Set-TextCell "E13" "-0.37%"
Set-TextCell "D14" "0.001260"
Set-TextCell "E14" "-2.60%"
Set-TextCell "D15" "0.005922"
Set-TextCell "E15" "4.97%"
Set-TextCell "D16" "3.356"
Set-TextCell "E16" "-0.10%"
Set-TextCell "D17" "2.436"
Set-TextCell "E17" "-2.72%"
Set-TextCell "E18" "-1.96%"
Set-TextCell "D19" "7.053"
Set-TextCell "E19" "-4.85%"
Set-TextCell "D20" "0.1338"
Set-TextCell "E20" "-3.65%"
Set-TextCell "D21" "0.3220"
Set-TextCell "E21" "22.84%"
Set-TextCell "D22" "0.04196"
Set-TextCell "E22" "0.38%"
Set-TextCell "E23" "-4.53%"
Set-TextCell "D24" "0.004063"
Set-TextCell "E24" "-5.40%"
Set-TextCell "D25" "0.0001220"
Set-TextCell "E25" "9.91%"
Set-TextCell "E26" "-0.01%"
Set-TextCell "D38" "0.02337"
Set-TextCell "E38" "-13.42%"
Set-TextCell "D39" "0.05146"
Set-TextCell "E39" "-6.61%"
Set-TextCell "D40" "0.007950"
Set-TextCell "E40" "0.03%"
Set-TextCell "D41" "0.1322"
Set-TextCell "E41" "-5.00%"
Set-TextCell "D42" "0.006377"
Set-TextCell "E42" "-2.70%"
Set-TextCell "D43" "0.001949"
Set-TextCell "E43" "-8.06%"
Set-TextCell "D44" "0.008594"
Set-TextCell "E44" "15.14%"
Set-TextCell "D45" "0.3044"
Set-TextCell "E45" "-5.01%"
Set-TextCell "D46" "0.00006511"
Set-TextCell "E46" "-6.83%"
Set-TextCell "E47" "-0.04%"
Set-TextCell "D48" "0.006999"
Set-TextCell "E48" "98.08%"
Set-TextCell "D49" "0.003371"
Set-TextCell "E49" "-2.74%"
Set-TextCell "D50" "0.00002100"
Set-TextCell "E50" "-0.04%"
Set-TextCell "D51" "0.0002000"
Set-TextCell "E51" "-0.04%"
